$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Update existing cell values (rows are 1-indexed, row 1 is the header) ---

# Row 2: Study design / Diagnostic test accuracy
$t.Cell(2,3).Range.Text = "28"
$t.Cell(2,4).Range.Text = "21.2"

# Row 4: Study design / Surveillance
$t.Cell(4,3).Range.Text = "93"
$t.Cell(4,4).Range.Text = "70.5"

# Row 5: Study design / Vaccine effectiveness
$t.Cell(5,3).Range.Text = "10"
$t.Cell(5,4).Range.Text = "7.6"

# Row 6: Sampling quality / High
$t.Cell(6,4).Range.Text = "28.0"

# Row 7: Sampling quality / Low
$t.Cell(7,3).Range.Text = "95"
$t.Cell(7,4).Range.Text = "72.0"

# Row 8: Prop lab tested / 0-5
$t.Cell(8,3).Range.Text = "12"
$t.Cell(8,4).Range.Text = "9.1"

# Row 9: Prop lab tested / 5-50
$t.Cell(9,3).Range.Text = "32"
$t.Cell(9,4).Range.Text = "24.2"

# Row 10: Prop lab tested / 50-95
$t.Cell(10,3).Range.Text = "27"
$t.Cell(10,4).Range.Text = "20.5"

# Row 11: Prop lab tested / 95+
$t.Cell(11,3).Range.Text = "30"
$t.Cell(11,4).Range.Text = "22.7"

# Row 12: was Num tests / 1 -> Prop lab tested / Not reported
$t.Cell(12,1).Range.Text = "Prop lab tested"
$t.Cell(12,2).Range.Text = "Not reported"
$t.Cell(12,3).Range.Text = "31"
$t.Cell(12,4).Range.Text = "23.5"

# Row 13: Num tests / 2 -> Num tests / 1
$t.Cell(13,2).Range.Text = "1"
$t.Cell(13,3).Range.Text = "106"
$t.Cell(13,4).Range.Text = "80.3"

# Row 14: Num tests / 3+ -> Num tests / 2
$t.Cell(14,2).Range.Text = "2"
$t.Cell(14,3).Range.Text = "19"
$t.Cell(14,4).Range.Text = "14.4"

# Row 15: was N / 1-9 -> Num tests / 3+
$t.Cell(15,1).Range.Text = "Num tests"
$t.Cell(15,2).Range.Text = "3+"
$t.Cell(15,3).Range.Text = "7"
$t.Cell(15,4).Range.Text = "5.3"

# Row 16: N / 10-99 -> N / 1-9
$t.Cell(16,2).Range.Text = "1-9"
$t.Cell(16,3).Range.Text = "1"
$t.Cell(16,4).Range.Text = "0.8"

# Row 17: N / 100-999 -> N / 10-99
$t.Cell(17,2).Range.Text = "10-99"
$t.Cell(17,3).Range.Text = "37"
$t.Cell(17,4).Range.Text = "28.0"

# --- Insert a brand-new row for "N / 100-999" right before the final row (old row 18, N / 1000+) ---
$newRow = $t.Rows.Add($t.Rows(18))
$t.Cell(18,1).Range.Text = "N"
$t.Cell(18,2).Range.Text = "100-999"
$t.Cell(18,3).Range.Text = "55"
$t.Cell(18,4).Range.Text = "41.7"

# --- Final row (now row 19, was row 18: N / 1000+) ---
$t.Cell(19,3).Range.Text = "39"
$t.Cell(19,4).Range.Text = "29.5"
